$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: the "o" marker in E3 was removed by the user -> clear the cell
$ws.Range("E3").ClearContents()

# Row 5: A5 used to hold the licence number "35" as text; store it as a
# real number instead (matches how row 2/3/4 already store it).
$ws.Range("A5").Value = 35

# New row 6: another competitor result row. The licence number must stay
# text rather than become a number, so format the cell as Text first
# (same effect as a person typing it in manually), then drop back to the
# Normal style so no extra formatting lingers on the cell.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "6858"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "Fekete Kálmán"
$ws.Range("C6").Value = "Üllői Lövész Klub"
$ws.Range("V6").Value = "VID_00001"

# New row 7: a mostly blank row, only the competition id column is filled.
$ws.Range("V7").Value = "VID_00001"
